# Applies the "Updated cryptos list" data refresh to Sheet1 (Price + Volume(1h) columns).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.466.46"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").Value = "1.884.13"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'0.7203"
$ws.Range("E5").Value = "  +1.79%  "
$ws.Range("D6").Value = "'243.77"
$ws.Range("E6").Value = "  +0.81%  "
$ws.Range("D8").Value = "'0.07977"
$ws.Range("E8").Value = "  +2.29%  "
$ws.Range("D9").Value = "'0.3156"
$ws.Range("E9").Value = "  +1.57%  "
$ws.Range("D10").Value = "'25.09"
$ws.Range("E10").Value = "  +0.26%  "
$ws.Range("D11").Value = "'0.08147"
$ws.Range("E11").Value = "  -2.96%  "
$ws.Range("D12").Value = "1.896.01"
$ws.Range("E12").Value = "  +1.17%  "
$ws.Range("D13").Value = "'5.263"
$ws.Range("E13").Value = "  +0.51%  "
$ws.Range("D14").Value = "'94.96"
$ws.Range("E14").Value = "  +4.30%  "
$ws.Range("D15").Value = "'0.7124"
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("D16").Value = "'6.408"
$ws.Range("E16").Value = "  +4.48%  "
$ws.Range("D17").Value = "'0.000008459"
$ws.Range("E17").Value = "  +1.26%  "
$ws.Range("D18").Value = "29.465.48"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("D19").Value = "'254.71"
$ws.Range("E19").Value = "  +6.03%  "
$ws.Range("E20").Value = "  +1.19%  "
$ws.Range("D21").Value = "2.141.83"
$ws.Range("E21").Value = "  +0.97%  "
$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").Value = "'7.807"
$ws.Range("E23").Value = "  +0.73%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").Value = "'0.1593"
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("D26").Value = "'9.095"
$ws.Range("E26").Value = "  +0.73%  "
$ws.Range("D27").Value = "'162.85"
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("D28").Value = "'19.07"
$ws.Range("E28").Value = "  +3.06%  "
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("D30").Value = "'4.432"
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("D31").Value = "'4.299"
$ws.Range("E31").Value = "  -0.66%  "
$ws.Range("E32").Value = "  -2.42%  "
$ws.Range("D33").Value = "'0.05333"
$ws.Range("E33").Value = "  -0.35%  "
$ws.Range("D34").Value = "'1.959"
$ws.Range("E34").Value = "  +1.12%  "
$ws.Range("E35").Value = "  +1.15%  "
$ws.Range("D36").Value = "'1.185"
$ws.Range("E36").Value = "  +0.89%  "
$ws.Range("D37").Value = "'2.701"
$ws.Range("E37").Value = "  +0.68%  "
$ws.Range("E38").Value = "  +1.25%  "
$ws.Range("D39").Value = "1.276.52"
$ws.Range("E39").Value = "  +2.97%  "
$ws.Range("D40").Value = "'2.770"
$ws.Range("E40").Value = "  +1.34%  "
$ws.Range("D41").Value = "'6.480"
$ws.Range("E41").Value = "  -0.51%  "
$ws.Range("D42").Value = "'113.18"
$ws.Range("E42").Value = "  +3.70%  "
$ws.Range("D43").Value = "'74.63"
$ws.Range("E43").Value = "  +3.16%  "
$ws.Range("D44").Value = "'0.9056"
$ws.Range("E44").Value = "  +1.71%  "
$ws.Range("D45").Value = "'0.00000000130"
$ws.Range("E45").Value = "  +3.10%  "
$ws.Range("D47").Value = "2.038.38"
$ws.Range("E47").Value = "  +0.86%  "
$ws.Range("D48").Value = "'1.807"
$ws.Range("E48").Value = "  +1.01%  "
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("D50").Value = "'9.520"
$ws.Range("E50").Value = "  +0.82%  "
$ws.Range("D51").Value = "'0.4385"
$ws.Range("E51").Value = "  +1.03%  "
